$d = $word.ActiveDocument

$pkgHeader = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# 1. Paragraph 2: "Checkout: removed disabled pa with amazon feature."
#    -> "Footer: Corrected typo in footer link." (numId 20 -> 22, drop underline rPr, drop bookmark)
$p2Xml = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="22"/></w:numPr></w:pPr><w:r><w:t>Footer: Corrected typo in footer link.</w:t></w:r></w:p>'
$d.Paragraphs.Item(2).Range.InsertXML($pkgHeader + $p2Xml + $pkgFooter)

# 2. Paragraph 3: "Checkout: Corrected tiny text in the cart page for cart totals."
#    -> "Product: Completed reviews section- no reviews message if no reviews are present." (numId 20 -> 22, drop underline rPr)
$p3Xml = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="22"/></w:numPr></w:pPr><w:r><w:t>Product: Completed reviews section- no reviews message if no reviews are present.</w:t></w:r></w:p>'
$d.Paragraphs.Item(3).Range.InsertXML($pkgHeader + $p3Xml + $pkgFooter)

# 3. Delete paragraph 4 entirely: "Header: Added truck image."
$d.Paragraphs.Item(4).Range.Delete()

# 4. Delete the block of numId=19 bullet paragraphs plus the trailing empty paragraph
#    (originally paragraphs 6-10 [after "To Do:"], now paragraphs 5-9 after the previous deletion)
$blockStart = $d.Paragraphs.Item(5)
$blockEnd = $d.Paragraphs.Item(9)
$d.Range($blockStart.Range.Start, $blockEnd.Range.End).Delete()

# 5. Simplify the "Badges" paragraph runs (remove proofErr spell-check wrappers, merge trailing runs)
#    This paragraph is now at index 6.
$badgesXml = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="18"/></w:numPr></w:pPr><w:r><w:t>Badges, where do we get those from?</w:t></w:r><w:r><w:t xml:space="preserve"> Specifically: google checkout and mcafee secure.</w:t></w:r></w:p>'
$d.Paragraphs.Item(6).Range.InsertXML($pkgHeader + $badgesXml + $pkgFooter)

# 6. Remove the now-orphaned _GoBack bookmark if it still exists
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
